$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.5921199760866448
$ws.Range("E2").Value = 0.5921199760866448

# Row 3
$ws.Range("D3").Value = 0.002613423715615775
$ws.Range("E3").Value = 0.002613423715615775

# Row 4
$ws.Range("D4").Value = [double]"2.736856800357776E-05"
$ws.Range("E4").Value = [double]"2.736856800357776E-05"

# Row 5
$ws.Range("D5").Value = 0.06081060715217985
$ws.Range("E5").Value = 0.06081060715217985

# Row 6
$ws.Range("D6").Value = 0.935004063684319
$ws.Range("E6").Value = 0.935004063684319

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.3548207819466465
$ws.Range("E7").Value = 0.6451792180533535

# Row 8
$ws.Range("D8").Value = 0.9999769332025726
$ws.Range("E8").Value = [double]"2.30667974273624E-05"

# Row 9
$ws.Range("D9").Value = 0.9943894140727879
$ws.Range("E9").Value = 0.005610585927212086

# Row 10
$ws.Range("D10").Value = 0.9999990182449792
$ws.Range("E10").Value = [double]"9.817550208346404E-07"

# Row 11
$ws.Range("D11").Value = 0.999998364899359
$ws.Range("E11").Value = [double]"1.635100640995368E-06"
$ws.Range("F11").Value = 0.4737389087677002
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.9935451281299104
$ws.Range("E12").Value = 0.9935451281299104

# Row 13
$ws.Range("D13").Value = 0.000107345042758357
$ws.Range("E13").Value = 0.000107345042758357

# Row 14
$ws.Range("D14").Value = [double]"5.023530983791071E-08"
$ws.Range("E14").Value = [double]"5.023530983791071E-08"

# Row 15
$ws.Range("D15").Value = 0.007650722287259481
$ws.Range("E15").Value = 0.007650722287259481

# Row 16
$ws.Range("D16").Value = 0.9788794552197997
$ws.Range("E16").Value = 0.9788794552197997

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.01340822701728257
$ws.Range("E17").Value = 0.9865917729827174

# Row 18
$ws.Range("D18").Value = 0.9999999992921746
$ws.Range("E18").Value = [double]"7.078253538850277E-10"

# Row 19
$ws.Range("D19").Value = 0.9985298249144169
$ws.Range("E19").Value = 0.001470175085583114

# Row 20
$ws.Range("D20").Value = 0.9999968509380414
$ws.Range("E20").Value = [double]"3.149061958640509E-06"

# Row 21
$ws.Range("D21").Value = 0.999994889272796
$ws.Range("E21").Value = [double]"5.110727203994081E-06"
$ws.Range("F21").Value = 1.322158455848694
$ws.Range("G21").Value = 0.7
